$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138, shifting existing rows 138-196 down to 139-197.
$ws.Rows.Item(138).Insert()

# Populate the newly inserted row 138 with its data (columns A-K mirror the
# surrounding rows for this market/product block; D, L, M, N, O, P, Q, R, S, T
# are the new record's own values).
$ws.Cells.Item(138, 1).Value = 3
$ws.Cells.Item(138, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(138, 3).Value = "Coquimbo"
$ws.Cells.Item(138, 4).Value = 44609
$ws.Cells.Item(138, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(138, 5).Value = 5
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100101
$ws.Cells.Item(138, 8).Value = "Berries"
$ws.Cells.Item(138, 9).Value = 100101001
$ws.Cells.Item(138, 10).Value = "Arándano (blue)"
$ws.Cells.Item(138, 11).Value = "Sin especificar"
$ws.Cells.Item(138, 12).Value = "Primera"
$ws.Cells.Item(138, 13).Value = 30
$ws.Cells.Item(138, 14).Value = 4000
$ws.Cells.Item(138, 15).Value = 4000
$ws.Cells.Item(138, 16).Value = 4000
$ws.Cells.Item(138, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(138, 18).Value = "Provincia de Linares"
$ws.Cells.Item(138, 19).Value = 2000
$ws.Cells.Item(138, 20).Value = 2
